$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of "year data" (A column index continues 204..215, B column values)
$data = @(
    @(204, 0.5217391304347826),
    @(205, 0.7646424721523535),
    @(206, 0.3162055335968379),
    @(207, 0.3162055335968379),
    @(208, 0.5059288537549407),
    @(209, 0.5059288537549407),
    @(210, 0.7905138339920948),
    @(211, 0.4869565217391304),
    @(212, 0.5059288537549407),
    @(213, 0.2529644268774703),
    @(214, 0.5059288537549407),
    @(215, 0.5059288537549407)
)

# Existing data ends at row 205 (A205 = 203, B205 = last value). New data starts at row 206.
$startRow = 206

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $aVal = $data[$i][0]
    $bVal = $data[$i][1]

    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)

    $cellA.Value = $aVal
    $cellB.Value = $bVal
}

# Match formatting of column A in the pre-existing data block (centered/bold/bordered style)
$srcA = $ws.Range("A205")
$dstA = $ws.Range("A206:A217")
$srcA.Copy()
$dstA.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the sheet dimension to reflect the new used range
$ws.UsedRange | Out-Null
